$d = $word.ActiveDocument
$sec = $d.Sections(1)

# Helper: rename the single inline picture living in a HeaderFooter's range.
# InlineShape has no writable .Name in the Word object model, so the
# picture has to be converted to a floating Shape (which does expose
# .Name), renamed, and converted straight back to an inline shape so the
# surrounding <w:drawing><wp:inline> layout is preserved.
function Rename-InlinePicture($range, $newName) {
    $shape = $range.InlineShapes(1)
    $floatShape = $shape.ConvertToShape()
    $floatShape.Name = $newName
    $floatShape.ConvertToInlineShape() | Out-Null
}

# Footer (default/primary, physically footer2.xml) - Pearson logo: image2.png -> image1.png
Rename-InlinePicture $sec.Footers(1).Range "image1.png"

# Footer (first page, physically footer1.xml) - Pearson logo: image2.png -> image1.png
Rename-InlinePicture $sec.Footers(2).Range "image1.png"

# Header (first page, physically header1.xml) - BTEC logo: image1.jpg -> image2.jpg
Rename-InlinePicture $sec.Headers(2).Range "image2.jpg"
